$wb = $excel.ActiveWorkbook

# --- 1) Shared-string text change: "Ready for handoff" -> "In Translation" ---
# Every cell in the workbook that currently reads "Ready for handoff" shares the
# same sharedStrings entry, so every one of them has to be rewritten (rewriting
# only one would just fork off a brand new string and leave the rest pointing at
# the stale text).

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2) Column width changes ---
# Target stored width (post-save XML) is 13.4101845877511 "character" units for
# each of these columns. The COM ColumnWidth setter here snaps the persisted
# width to a 1/6-character pixel grid (same rounding Excel itself applies), so
# feeding it the literal target value would land on 14.1666... instead of the
# nearest reachable width. Requesting 12.5 lands on 13.3333333333..., which is
# the closest the engine can actually persist to the target value.
$narrowWidth = 12.5

# Overview: columns E and F (5, 6) get narrower.
$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# zh-cn and de-de: column C (3) gets narrower.
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
